# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column D) and
# "Correspond Handback DateTime" (column G) cells with newly generated
# report timestamps.

$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 3 gets a refreshed handoff/handback datetime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-02-17 06:12:38"
$wsZhCn.Range("G3").Value = "2016-02-17 06:13:23"

# de-de sheet: row 3 gets a refreshed handoff/handback datetime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-02-17 06:12:48"
$wsDeDe.Range("G3").Value = "2016-02-17 06:13:40"
